# "filled in some more bare bones" -- add body text to several slides
# that only had empty placeholder paragraphs.

$p = $ppt.ActivePresentation

# --- Slide 3 ("The Old" / "Does these Things well") ---------------------
# Content placeholder just gets a second, still-empty paragraph.
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "`r"

# --- Slide 4 ("The Old" / "Does these Things poorly") --------------------
$s4 = $p.Slides.Item(4)
$cp4 = $s4.Shapes.Item(1)
$cp4.TextFrame.TextRange.Text = "Routing Problems`rCount to infinity`rNetwork Topology"
$cp4.TextFrame.TextRange.Paragraphs(2).IndentLevel = 2

# --- Slide 5 ("The Old" / "What's the flow") ------------------------------
$s5 = $p.Slides.Item(5)
$cp5 = $s5.Shapes.Item(1)
$cp5.Left = 96
$cp5.Top = 204
$cp5.Width = 522
$cp5.Height = 228
$cp5.TextFrame.TextRange.Text = "Graphic of networking flow goes here"
$cp5.TextFrame.TextRange.Font.Size = 60

# --- Slide 7 ("The New" / "What's Does it fix?") --------------------------
$s7 = $p.Slides.Item(7)
$cp7 = $s7.Shapes.Item(1)
$cp7.TextFrame.TextRange.Text = "Routing`rNetwork Topology"

# --- Slide 8 ("The New" / "Problems") -------------------------------------
$s8 = $p.Slides.Item(8)
$cp8 = $s8.Shapes.Item(1)
$cp8.TextFrame.TextRange.Text = "Single point of failure`rFinding a viable, scalable, universal protocol for all switches/ controllers to talk through (currently OpenFlow, but it is weird)"

# --- Slide 9 ("The New" / "Flow") -----------------------------------------
# The original content placeholder is removed and a fresh one (re-sized,
# moved to the end of the shape list, with the graphic caption) takes
# its place.
$s9 = $p.Slides.Item(9)
$origCp9 = $s9.Shapes.Item(1)
$newCp9 = $origCp9.Duplicate()
$origCp9.Delete()
$newCp9.Left = 96
$newCp9.Top = 204
$newCp9.Width = 522
$newCp9.Height = 228
$newCp9.TextFrame.TextRange.Text = "Graphic of networking flow goes here"
$newCp9.TextFrame.TextRange.Font.Size = 60
